# Rename the inline picture shapes so their Word "Name" (the OOXML
# wp:docPr/@name - and, via the picture's non-visual properties,
# pic:cNvPr/@name) matches the target filenames from the diff:
#
#   - Footer "first page"  (footer1.xml, docPr id=3): image1.png -> image2.png
#   - Footer "default"     (footer2.xml, docPr id=2): image1.png -> image2.png
#   - Header "first page"  (header1.xml, docPr id=1): image2.jpg -> image1.jpg
#
# Both the Pearson logo footers and the BTec logo header only have a
# "first page" + "default/primary" header/footer pair (no distinct even
# page header/footer), and there is a single section in this document.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2

# --- Footers: Pearson Edexcel logo (image1.png -> image2.png) ---
$footerPrimary = $sec.Footers.Item(1)
$pearsonPrimary = $footerPrimary.Range.InlineShapes.Item(1)
$pearsonPrimary.Name = "image2.png"

$footerFirst = $sec.Footers.Item(2)
$pearsonFirst = $footerFirst.Range.InlineShapes.Item(1)
$pearsonFirst.Name = "image2.png"

# --- Header: BTec logo (image2.jpg -> image1.jpg) ---
$headerFirst = $sec.Headers.Item(2)
$btec = $headerFirst.Range.InlineShapes.Item(1)
$btec.Name = "image1.jpg"
